$wb = $excel.ActiveWorkbook

# Sheet ALC (sheet1.xml), row 15
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1110.3334
$ws.Range("I15").Value = 1110.3334
$ws.Range("K15").Value = 3331.0002
$ws.Range("M15").Value = -3162.0002

# Sheet ALC (sheet1.xml), row 130
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H130").Value = 48830
$ws.Range("J130").Value = 48830
$ws.Range("L130").Value = 48830
$ws.Range("N130").Value = -58870

# Sheet ARM (sheet2.xml), row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 27048
$ws.Range("I32").Value = 27689.947
$ws.Range("J32").Value = 23053.666
$ws.Range("K32").Value = 27689.947
$ws.Range("L32").Value = 23053.666
$ws.Range("M32").Value = -27402.947
$ws.Range("N32").Value = -23627.666

# Sheet ARM (sheet2.xml), row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1279.0714
$ws.Range("I74").Value = 653.1429000000001
$ws.Range("J74").Value = 3156.8572
$ws.Range("K74").Value = 653.1429000000001
$ws.Range("L74").Value = 3156.8572
$ws.Range("M74").Value = 220.8570999999999
$ws.Range("N74").Value = -4904.8572

# Sheet ARM (sheet2.xml), row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 1279.0714
$ws.Range("I77").Value = 653.1429000000001
$ws.Range("J77").Value = 3156.8572
$ws.Range("K77").Value = 3265.7145
$ws.Range("L77").Value = 15784.286
$ws.Range("M77").Value = 1102.2855
$ws.Range("N77").Value = -24520.286

# Sheet ARM (sheet2.xml), row 133
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H133").Value = 33462.582
$ws.Range("J133").Value = 33462.582
$ws.Range("L133").Value = 33462.582
$ws.Range("N133").Value = -38522.582

# Sheet ARM (sheet2.xml), row 134
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H134").Value = 50788
$ws.Range("J134").Value = 50788
$ws.Range("L134").Value = 50788
$ws.Range("N134").Value = -60928

# Sheet BSM (sheet3.xml), row 82
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 2503
$ws.Range("I82").Value = 2503
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 2503
$ws.Range("L82").Value = 0
$ws.Range("M82").Value = -2120
$ws.Range("N82").ClearContents()

# Sheet BSM (sheet3.xml), row 85
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H85").Value = 2503
$ws.Range("I85").Value = 2503
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 2503
$ws.Range("L85").Value = 0
$ws.Range("M85").Value = -1177
$ws.Range("N85").ClearContents()

# Sheet BSM (sheet3.xml), row 117
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H117").Value = 48933.5
$ws.Range("J117").Value = 48933.5
$ws.Range("L117").Value = 48933.5
$ws.Range("N117").Value = -58111.5

# Sheet BSM (sheet3.xml), row 130
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H130").Value = 38492.855
$ws.Range("J130").Value = 38492.855
$ws.Range("L130").Value = 38492.855
$ws.Range("N130").Value = -48532.855

# Sheet BSM (sheet3.xml), row 132
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H132").Value = 41795
$ws.Range("J132").Value = 41795
$ws.Range("L132").Value = 41795
$ws.Range("N132").Value = -51915

# Sheet BSM (sheet3.xml), row 137
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H137").Value = 28734.834
$ws.Range("J137").Value = 28734.834
$ws.Range("L137").Value = 28734.834
$ws.Range("N137").Value = -38934.834

# Sheet CRP (sheet4.xml), row 58
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1620.2703
$ws.Range("I58").Value = 1362.1613
$ws.Range("J58").Value = 2953.8333
$ws.Range("K58").Value = 1362.1613
$ws.Range("L58").Value = 2953.8333
$ws.Range("M58").Value = -1159.1613
$ws.Range("N58").Value = -3359.8333

# Sheet CRP (sheet4.xml), row 115
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H115").Value = 34255
$ws.Range("J115").Value = 34255
$ws.Range("L115").Value = 34255
$ws.Range("N115").Value = -36605

# Sheet CRP (sheet4.xml), row 118
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H118").Value = 43734
$ws.Range("J118").Value = 43734
$ws.Range("L118").Value = 43734
$ws.Range("N118").Value = -47048

# Sheet CRP (sheet4.xml), row 136
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 1620.2703
$ws.Range("I136").Value = 1362.1613
$ws.Range("J136").Value = 2953.8333
$ws.Range("K136").Value = 4086.4839
$ws.Range("L136").Value = 8861.499899999999
$ws.Range("M136").Value = -1536.4839
$ws.Range("N136").Value = -13961.4999

# Sheet CUL (sheet5.xml), row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 22007.21
$ws.Range("I131").Value = 10524.9
$ws.Range("J131").Value = 23283.021
$ws.Range("K131").Value = 31574.7
$ws.Range("L131").Value = 69849.06299999999
$ws.Range("M131").Value = -26534.7
$ws.Range("N131").Value = -79929.06299999999

# Sheet GSM (sheet6.xml), row 110
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H110").Value = 28075.2
$ws.Range("J110").Value = 28075.2
$ws.Range("L110").Value = 28075.2
$ws.Range("N110").Value = -36255.2

# Sheet GSM (sheet6.xml), row 135
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H135").Value = 32127.363
$ws.Range("J135").Value = 32127.363
$ws.Range("L135").Value = 32127.363
$ws.Range("N135").Value = -42267.363

# Sheet GSM (sheet6.xml), row 138
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H138").Value = 51750
$ws.Range("J138").Value = 51750
$ws.Range("L138").Value = 51750
$ws.Range("N138").Value = -62030

# Sheet LTW (sheet7.xml), row 63
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H63").Value = 20085
$ws.Range("J63").Value = 20085
$ws.Range("L63").Value = 20085
$ws.Range("N63").Value = -21583

# Sheet LTW (sheet7.xml), row 66
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H66").Value = 20085
$ws.Range("J66").Value = 20085
$ws.Range("L66").Value = 60255
$ws.Range("N66").Value = -67743

# Sheet LTW (sheet7.xml), row 111
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H111").Value = 43934.5
$ws.Range("J111").Value = 43934.5
$ws.Range("L111").Value = 43934.5
$ws.Range("N111").Value = -52114.5

# Sheet LTW (sheet7.xml), row 121
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H121").Value = 18891.334
$ws.Range("J121").Value = 18891.334
$ws.Range("L121").Value = 18891.334
$ws.Range("N121").Value = -22385.334

# Sheet LTW (sheet7.xml), row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3987.52
$ws.Range("I132").Value = 2290.818
$ws.Range("J132").Value = 5320.643
$ws.Range("K132").Value = 6872.454000000001
$ws.Range("L132").Value = 15961.929
$ws.Range("M132").Value = -4342.454000000001
$ws.Range("N132").Value = -21021.929

# Sheet LTW (sheet7.xml), row 141
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H141").Value = 14648.167
$ws.Range("J141").Value = 14648.167
$ws.Range("L141").Value = 14648.167
$ws.Range("N141").Value = -25008.167

# Sheet WVR (sheet8.xml), row 119
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H119").Value = 45729.668
$ws.Range("J119").Value = 45729.668
$ws.Range("L119").Value = 45729.668
$ws.Range("N119").Value = -55405.668

# Sheet WVR (sheet8.xml), row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2084.879
$ws.Range("I132").Value = 1560.4
$ws.Range("J132").Value = 2891.7693
$ws.Range("K132").Value = 4681.200000000001
$ws.Range("L132").Value = 8675.3079
$ws.Range("M132").Value = -2151.200000000001
$ws.Range("N132").Value = -13735.3079

# Sheet WVR (sheet8.xml), row 133
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H133").Value = 82064.25
$ws.Range("J133").Value = 82064.25
$ws.Range("L133").Value = 82064.25
$ws.Range("N133").Value = -92184.25

# Sheet WVR (sheet8.xml), row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 16876.375
$ws.Range("I136").Value = 27000.736
$ws.Range("J136").Value = 2079.2307
$ws.Range("K136").Value = 81002.208
$ws.Range("L136").Value = 6237.6921
$ws.Range("M136").Value = -78452.208
$ws.Range("N136").Value = -11337.6921

# Sheet WVR (sheet8.xml), row 137
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H137").Value = 22550
$ws.Range("J137").Value = 22550
$ws.Range("L137").Value = 22550
$ws.Range("N137").Value = -32750
